$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 20150
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 20150
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 20150
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = -20374

$ws.Range("H14").Value = 20150
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 20150
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 20150
$ws.Range("M14").Value = ""
$ws.Range("N14").Value = -20532

$ws.Range("H70").Value = 1332.6666
$ws.Range("J70").Value = 1349.5
$ws.Range("L70").Value = 4048.5
$ws.Range("N70").Value = -4588.5

$ws.Range("H73").Value = 1332.6666
$ws.Range("J73").Value = 1349.5
$ws.Range("L73").Value = 4048.5
$ws.Range("N73").Value = -5920.5

$ws.Range("H74").Value = 4000
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 4000
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 4000
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 20000
$ws.Range("N77").Value = -29360

$ws.Range("H86").Value = 4996.3335
$ws.Range("I86").Value = 4995
$ws.Range("K86").Value = 4995
$ws.Range("M86").Value = -3872

$ws.Range("H89").Value = 4996.3335
$ws.Range("I89").Value = 4995
$ws.Range("K89").Value = 24975
$ws.Range("M89").Value = -19359

$ws.Range("H99").Value = 171.2
$ws.Range("I99").Value = 171.2
$ws.Range("K99").Value = 513.5999999999999
$ws.Range("M99").Value = 984.4000000000001

$ws.Range("H106").Value = 1199
$ws.Range("I106").Value = 1199
$ws.Range("K106").Value = 1199
$ws.Range("M106").Value = -568

$ws.Range("H115").Value = 4322.25
$ws.Range("I115").Value = 4322.25
$ws.Range("K115").Value = 12966.75
$ws.Range("M115").Value = -11399.75

$ws.Range("H118").Value = 1266.6666
$ws.Range("I118").Value = 1075
$ws.Range("K118").Value = 3225
$ws.Range("M118").Value = -1568

$ws.Range("H121").Value = 1751.8334
$ws.Range("J121").Value = 1902.4546
$ws.Range("L121").Value = 5707.3638
$ws.Range("N121").Value = -9201.363799999999

$ws.Range("H137").Value = 3302.6775
$ws.Range("I137").Value = 1970.762
$ws.Range("J137").Value = 6099.7
$ws.Range("K137").Value = 5912.286
$ws.Range("L137").Value = 18299.1
$ws.Range("M137").Value = -3362.286
$ws.Range("N137").Value = -23399.1

$ws.Range("H138").Value = 3013.2954
$ws.Range("I138").Value = 1262.1305
$ws.Range("J138").Value = 4931.2383
$ws.Range("K138").Value = 3786.3915
$ws.Range("L138").Value = 14793.7149
$ws.Range("M138").Value = 1353.6085
$ws.Range("N138").Value = -25073.7149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2741.4285
$ws.Range("I61").Value = 2738.6
$ws.Range("J61").Value = 2748.5
$ws.Range("K61").Value = 2738.6
$ws.Range("L61").Value = 2748.5
$ws.Range("M61").Value = -2526.6
$ws.Range("N61").Value = -3172.5

$ws.Range("H132").Value = 2267.182
$ws.Range("I132").Value = 1887.8148
$ws.Range("K132").Value = 5663.4444
$ws.Range("M132").Value = -3133.4444

$ws.Range("H136").Value = 2741.4285
$ws.Range("I136").Value = 2738.6
$ws.Range("J136").Value = 2748.5
$ws.Range("K136").Value = 8215.799999999999
$ws.Range("L136").Value = 8245.5
$ws.Range("M136").Value = -5665.799999999999
$ws.Range("N136").Value = -13345.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2923.5334
$ws.Range("I20").Value = 2991.4546
$ws.Range("J20").Value = 2736.75
$ws.Range("K20").Value = 2991.4546
$ws.Range("L20").Value = 2736.75
$ws.Range("M20").Value = -2744.4546
$ws.Range("N20").Value = -3230.75

$ws.Range("H105").Value = 3564.7144
$ws.Range("I105").Value = 3990.8
$ws.Range("K105").Value = 3990.8
$ws.Range("M105").Value = -2243.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6061.75
$ws.Range("I31").Value = 3374.25
$ws.Range("J31").Value = 8749.25
$ws.Range("K31").Value = 3374.25
$ws.Range("L31").Value = 8749.25
$ws.Range("M31").Value = -3079.25
$ws.Range("N31").Value = -9339.25

$ws.Range("H34").Value = 6061.75
$ws.Range("I34").Value = 3374.25
$ws.Range("J34").Value = 8749.25
$ws.Range("K34").Value = 3374.25
$ws.Range("L34").Value = 8749.25
$ws.Range("M34").Value = -3172.25
$ws.Range("N34").Value = -9153.25

$ws.Range("H58").Value = 2242.5
$ws.Range("I58").Value = 2090.8948
$ws.Range("J58").Value = 2818.6
$ws.Range("K58").Value = 2090.8948
$ws.Range("L58").Value = 2818.6
$ws.Range("M58").Value = -1887.8948
$ws.Range("N58").Value = -3224.6

$ws.Range("H132").Value = 2994.9111
$ws.Range("I132").Value = 2975.081
$ws.Range("J132").Value = 3086.625
$ws.Range("K132").Value = 8925.243
$ws.Range("L132").Value = 9259.875
$ws.Range("M132").Value = -6395.243
$ws.Range("N132").Value = -14319.875

$ws.Range("H134").Value = 4724.1333
$ws.Range("I134").Value = 5113.5835
$ws.Range("J134").Value = 3166.3333
$ws.Range("K134").Value = 15340.7505
$ws.Range("L134").Value = 9498.999899999999
$ws.Range("M134").Value = -12805.7505
$ws.Range("N134").Value = -14568.9999

$ws.Range("H136").Value = 2242.5
$ws.Range("I136").Value = 2090.8948
$ws.Range("J136").Value = 2818.6
$ws.Range("K136").Value = 6272.6844
$ws.Range("L136").Value = 8455.799999999999
$ws.Range("M136").Value = -3722.6844
$ws.Range("N136").Value = -13555.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 746.5
$ws.Range("J97").Value = 744
$ws.Range("L97").Value = 2232
$ws.Range("N97").Value = -3224

$ws.Range("H114").Value = 6755.5
$ws.Range("I114").Value = 3513.5
$ws.Range("J114").Value = 9997.5
$ws.Range("K114").Value = 10540.5
$ws.Range("L114").Value = 29992.5
$ws.Range("M114").Value = -7286.5
$ws.Range("N114").Value = -36500.5

$ws.Range("H117").Value = 1202.75
$ws.Range("J117").Value = 1440.6666
$ws.Range("L117").Value = 4321.9998
$ws.Range("N117").Value = -11205.9998

$ws.Range("H129").Value = 2178
$ws.Range("J129").Value = 2411
$ws.Range("L129").Value = 7233
$ws.Range("N129").Value = -17233

$ws.Range("H131").Value = 1298.909
$ws.Range("J131").Value = 1937
$ws.Range("L131").Value = 5811
$ws.Range("N131").Value = -15891

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 837829
$ws.Range("I3").Value = 3002.5
$ws.Range("J3").Value = 1255242.2
$ws.Range("K3").Value = 3002.5
$ws.Range("L3").Value = 1255242.2
$ws.Range("M3").Value = -2886.5
$ws.Range("N3").Value = -1255474.2

$ws.Range("H31").Value = 350
$ws.Range("I31").Value = 350
$ws.Range("K31").Value = 350
$ws.Range("M31").Value = -58

$ws.Range("H37").Value = 350
$ws.Range("I37").Value = 350
$ws.Range("K37").Value = 350
$ws.Range("M37").Value = -73

$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 15000
$ws.Range("K70").Value = 15000
$ws.Range("M70").Value = -14730

$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 15000
$ws.Range("K73").Value = 15000
$ws.Range("M73").Value = -14064

$ws.Range("H80").Value = 7665
$ws.Range("J80").Value = 10002.5
$ws.Range("L80").Value = 10002.5
$ws.Range("N80").Value = -11998.5

$ws.Range("H83").Value = 7665
$ws.Range("J83").Value = 10002.5
$ws.Range("L83").Value = 50012.5
$ws.Range("N83").Value = -59996.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 284
$ws.Range("I55").Value = 261.5
$ws.Range("J55").Value = 299
$ws.Range("K55").Value = 261.5
$ws.Range("L55").Value = 299
$ws.Range("M55").Value = -88.5
$ws.Range("N55").Value = -645

$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1200
$ws.Range("K82").Value = 1200
$ws.Range("M82").Value = -839

$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1200
$ws.Range("K85").Value = 1200
$ws.Range("M85").Value = 48

$ws.Range("H136").Value = 4301.2
$ws.Range("I136").Value = 4301.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12903.6
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10353.6
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2300
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

$ws.Range("H132").Value = 2077.6316
$ws.Range("I132").Value = 1598.5
$ws.Range("K132").Value = 4795.5
$ws.Range("M132").Value = -2265.5
